$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 11.91605408587611
$ws.Range("R2").Value = 107.244486772885
$ws.Range("S2").Value = 0.003389527972372358
$ws.Range("T2").Value = 0.003389527972372358
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 62.57172357811945
$ws.Range("R3").Value = 563.145512203075
$ws.Range("S3").Value = 0.01779856031360011
$ws.Range("T3").Value = 0.01779856031360011
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 31.1994113001811
$ws.Range("R4").Value = 280.7947017016299
$ws.Range("S4").Value = 0.008874689268896424
$ws.Range("T4").Value = 0.008874689268896424
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 74.08116395601333
$ws.Range("R5").Value = 666.7304756041199
$ws.Range("S5").Value = 0.02107242679877013
$ws.Range("T5").Value = 0.02107242679877013
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 166.9836491353971
$ws.Range("R6").Value = 1502.852842218574
$ws.Range("S6").Value = 0.0474985885087669
$ws.Range("T6").Value = 0.0474985885087669
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("S7").Value = 0.2494171752748373
$ws.Range("T7").Value = 0.2494171752748373
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 437.2077797091791
$ws.Range("R8").Value = 3934.870017382612
$ws.Range("S8").Value = 0.1243639873051246
$ws.Range("T8").Value = 0.1243639873051246
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 1038.124113940965
$ws.Range("R9").Value = 9343.117025468688
$ws.Range("S9").Value = 0.2952949607007814
$ws.Range("T9").Value = 0.2952949607007814
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 41.79539783315477
$ws.Range("R10").Value = 376.158580498393
$ws.Range("S10").Value = 0.01188872331821859
$ws.Range("T10").Value = 0.01188872331821859
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 219.4694704477262
$ws.Range("R11").Value = 1975.225234029535
$ws.Range("S11").Value = 0.06242820851627768
$ws.Range("T11").Value = 0.06242820851627768
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 109.4315113085038
$ws.Range("R12").Value = 984.883601776534
$ws.Range("S12").Value = 0.03112785205287058
$ws.Range("T12").Value = 0.03112785205287058
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 259.8386762237573
$ws.Range("R13").Value = 2338.548086013815
$ws.Range("S13").Value = 0.07391125073933193
$ws.Range("T13").Value = 0.07391125073933193
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 12.33523952009066
$ws.Range("R14").Value = 111.017155680816
$ws.Range("S14").Value = 0.003508765493840592
$ws.Range("T14").Value = 0.003508765493840592
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 64.77288471154667
$ws.Range("R15").Value = 582.95596240392
$ws.Range("S15").Value = 0.01842468177794405
$ws.Range("T15").Value = 0.01842468177794405
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 32.29695069357866
$ws.Range("R16").Value = 290.6725562422079
$ws.Range("S16").Value = 0.009186884937688411
$ws.Range("T16").Value = 0.009186884937688412
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 76.68720658188798
$ws.Range("R17").Value = 690.1848592369919
$ws.Range("S17").Value = 0.02181371702067898
$ws.Range("T17").Value = 0.02181371702067898
